# Insert one new data row before the current row 25 (shifts old rows
# 25..72 down to 26..73) and populate the new row 25 with the
# "Black Amber" / "Segunda" quality record described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 25 (and everything below it) down by one row.
$ws.Rows.Item(25).Insert()

# Fill in the newly-inserted row 25 with the full record.
$ws.Range("A25").Value = 7
$ws.Range("B25").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C25").Value = "Ñuble"
$ws.Range("D25").Value = 44952
$ws.Range("E25").Value = 16
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100103
$ws.Range("H25").Value = "Frutos de hueso (carozo)"
$ws.Range("I25").Value = 100103002
$ws.Range("J25").Value = "Ciruela"
$ws.Range("K25").Value = "Black Amber"
$ws.Range("L25").Value = "Segunda"
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 9000
$ws.Range("O25").Value = 9000
$ws.Range("P25").Value = 9000
$ws.Range("Q25").Value = '$/bandeja 18 kilos granel'
$ws.Range("R25").Value = "Provincia de Curicó"
$ws.Range("S25").Value = 500
$ws.Range("T25").Value = 18
